# edit.ps1 - applies the "Added one line for about me" commit to before.docx
#
# Summary of the change (see unified diff):
#   1. Drop the stray _GoBack bookmarkStart that originally sat before the
#      "RAJ BHARATH KANNAN" heading run (w:colFirst/w:colLast variant).
#   2. The phone-number paragraph ("9840859553") gains w:sz=26 on both the
#      paragraph mark run properties and the run itself, and becomes the new
#      home of the _GoBack bookmark: bookmarkStart now opens right after the
#      paragraph properties and bookmarkEnd closes immediately after the
#      paragraph (as a sibling of the <w:p> elements in the table cell).
#   3. The "ABOUT ME" paragraph is split: the sentence "...upcoming trends
#      and technologies." keeps its run, then five new runs spell out
#      "I love writing object oriented and adaptive code. "; a brand-new
#      Heading1-styled paragraph is started with "I" and the pre-existing
#      " bring in a lot of positive energy..." run becomes its second run.
#   4. The trailing empty cell at the bottom-right of the table loses its
#      w:ind w:left="342" paragraph indent.
#   5. The now-orphaned duplicate bookmarkEnd at the bottom of the table is
#      removed (the bookmark only needs the one pair created in step 2).
#   6. The section's top margin shrinks from 990 twips (49.5pt) to 630
#      twips (31.5pt) to make room for the extra line of text.
#
# Word's object model does not expose bookmark mutation (Bookmarks.Add /
# the raw bookmark tags) through the high level Paragraphs/Range API in
# this host, so the safest way to reproduce the restructuring exactly -
# including the bookmark move - is to hand Word the finished OOXML for the
# body and let Range.InsertXML (a real Word COM method) splice it in, the
# same way pasting "Keep Source Formatting" XML would.

$d = $word.ActiveDocument

$newBodyXml = '<w:tbl><w:tblPr><w:tblW w:w="10800" w:type="dxa"/><w:tblInd w:w="-720" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="7380"/><w:gridCol w:w="3420"/></w:tblGrid><w:tr w:rsidR="00485411" w:rsidRPr="00A22C3A" w:rsidTr="00E779E0"><w:tc><w:tcPr><w:tcW w:w="7380" w:type="dxa"/></w:tcPr><w:p w:rsidR="00485411" w:rsidRPr="0013737E" w:rsidRDefault="0013737E" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="62"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="0013737E"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="62"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>RAJ BHARATH KANNAN</w:t></w:r></w:p><w:p w:rsidR="0013737E" w:rsidRPr="0013737E" w:rsidRDefault="0013737E" w:rsidP="0013737E"><w:pPr><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="0013737E"><w:rPr><w:sz w:val="24"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>SYSTEM ENGINEER, TATA CONSULTANCY SERVICES LIMITED</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3420" w:type="dxa"/></w:tcPr><w:p w:rsidR="00485411" w:rsidRPr="00A259A1" w:rsidRDefault="0013737E" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r w:rsidRPr="00A259A1"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:t>9840859553</w:t></w:r></w:p><w:bookmarkEnd w:id="0"/><w:p w:rsidR="0013737E" w:rsidRPr="000A26B4" w:rsidRDefault="000A26B4" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "mailto:rbrajbharath1@gmail.com" </w:instrText></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00EF759C" w:rsidRPr="000A26B4"><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:t>Mail Me</w:t></w:r></w:p><w:p w:rsidR="0013737E" w:rsidRPr="0013737E" w:rsidRDefault="000A26B4" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:hyperlink r:id="rId6" w:history="1"><w:r w:rsidR="008C1B3A" w:rsidRPr="000A26B4"><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:t>My Repo</w:t></w:r></w:hyperlink></w:p><w:p w:rsidR="0013737E" w:rsidRPr="0013737E" w:rsidRDefault="00866CED" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:hyperlink r:id="rId7" w:history="1"><w:r w:rsidR="0013737E" w:rsidRPr="00866CED"><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Linked </w:t></w:r><w:r w:rsidR="008C1B3A" w:rsidRPr="00866CED"><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:t>I</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00866CED"><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:t>n</w:t></w:r></w:hyperlink></w:p><w:p w:rsidR="0013737E" w:rsidRDefault="00866CED" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:hyperlink r:id="rId8" w:history="1"><w:r w:rsidR="0013737E" w:rsidRPr="00866CED"><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="26"/><w:szCs w:val="28"/></w:rPr><w:t>Twitter</w:t></w:r></w:hyperlink></w:p><w:p w:rsidR="000B2CA0" w:rsidRPr="0013737E" w:rsidRDefault="000B2CA0" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="0013737E" w:rsidRPr="00A22C3A" w:rsidTr="00E779E0"><w:tc><w:tcPr><w:tcW w:w="7380" w:type="dxa"/></w:tcPr><w:p w:rsidR="0013737E" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="0013737E"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>ABOUT ME</w:t></w:r></w:p><w:p w:rsidR="0013737E" w:rsidRDefault="00627399" w:rsidP="0013737E"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>I’m</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00EF759C"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>a passionate software engineer</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> with an unquenching thirst for learning. </w:t></w:r><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>This</w:t></w:r><w:r w:rsidR="00EF759C"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> passion drives</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> me to keep myself updated with upc</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">oming trends and technologies. </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">I love writing object </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">oriented and </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>adaptive</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">code. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>I</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> bring in a lot of positive energy into the team and create a lively work environment.</w:t></w:r></w:p><w:p w:rsidR="00001D44" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3420" w:type="dxa"/></w:tcPr><w:p w:rsidR="0013737E" w:rsidRPr="00627399" w:rsidRDefault="0013737E" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>SKILLS</w:t></w:r></w:p><w:p w:rsidR="0013737E" w:rsidRPr="00627399" w:rsidRDefault="0013737E" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Android Application Development</w:t></w:r></w:p><w:p w:rsidR="0013737E" w:rsidRPr="00A22C3A" w:rsidRDefault="0013737E" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Web Application Development</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00485411" w:rsidRPr="00A22C3A" w:rsidTr="00E779E0"><w:tc><w:tcPr><w:tcW w:w="7380" w:type="dxa"/></w:tcPr><w:p w:rsidR="00001D44" w:rsidRDefault="00001D44" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>PROFESSIONAL SUMMARY</w:t></w:r></w:p><w:p w:rsidR="00A22C3A" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>SYSTEMS ENGINEER (DEC 2013 ONWARDS)</w:t></w:r></w:p><w:p w:rsidR="00485411" w:rsidRDefault="00001D44" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="000B2CA0"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>TATA CONSULTANCY SERVICES LIMITED</w:t></w:r></w:p><w:p w:rsidR="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>W</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">orked for one of the world’s leading </w:t></w:r><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>banks,</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> involved with requirement analysis, test scenario identification, defect management</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00883626" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>D</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>eveloped a knowledge transfer tool for effective training of the new joiner associates</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00883626" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">howcased a </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>POC</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> on mobile automation testing for mobile web and app</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00883626" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>T</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>rained fellow team mates on tool development</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00883626" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>ASSISTANT SYSTEMS ENGINEER (DEC 2011 – DEC 2013)</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>TATA CONSULTANCY SERVICES LIMITED</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>D</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">eveloped a </w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>test</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="000B2CA0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>script scheduling tool</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> called 1 click scheduler that saved hours of manual efforts </w:t></w:r></w:p><w:p w:rsidR="0013737E" w:rsidRPr="00001D44" w:rsidRDefault="0013737E" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>W</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>orked on existing automation framework by adding the auto defect function to automatically log the defects in database</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00883626" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="000B2CA0" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>W</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>orked for one of the world’s leading banks as a</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>n assistant</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> systems</w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> engineer</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3420" w:type="dxa"/></w:tcPr><w:p w:rsidR="00485411" w:rsidRPr="00627399" w:rsidRDefault="00627399" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>PROGRAMMING LANGUAGES</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00627399" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>A</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>ndroid</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00627399" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>J</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>ava</w:t></w:r></w:p><w:p w:rsidR="00485411" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>C</w:t></w:r></w:p><w:p w:rsidR="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:ind w:left="342"/><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:ind w:left="342"/><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="004F17C4" w:rsidRPr="00627399" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>OTHER INTERESTS</w:t></w:r></w:p><w:p w:rsidR="004F17C4" w:rsidRPr="00627399" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Developing Pair Programmed Android applications</w:t></w:r></w:p><w:p w:rsidR="004F17C4" w:rsidRPr="00627399" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="004F17C4" w:rsidRPr="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:ind w:left="342"/><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Sharing my knowledge with peer associates on Android, Hadoop and other known technologies</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00485411" w:rsidRPr="00A22C3A" w:rsidTr="00E779E0"><w:tc><w:tcPr><w:tcW w:w="7380" w:type="dxa"/></w:tcPr><w:p w:rsidR="00001D44" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>PERSONAL PROJECTS</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:caps/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>MEETING HALL APPLICATION</w:t></w:r></w:p><w:p w:rsidR="00485411" w:rsidRPr="00001D44" w:rsidRDefault="00582001" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>D</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">eveloped a meeting scheduler android application which helps book a meeting room online and eliminates the availability of it to others in that time frame. </w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00883626" w:rsidP="00A22C3A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">SHAKE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>SHAKE</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>D</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>eveloped an android application for professionals att</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>ending long duration meetings. I</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>n case an important discussion is interrupted by an unnecessary call, just shake the phone to turn it to silent mode or to reject the call.</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00883626" w:rsidP="00A22C3A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>TAMIL PULI (EXCLUSIVE FOR TABLETS)</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>D</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve">eveloped an android application (exclusively for tablets) for children to read, learn and write </w:t></w:r><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Tamil</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> alphabets in a friendly animated environment. </w:t></w:r><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>This</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> application was designed for an </w:t></w:r><w:r w:rsidR="0013737E" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>NGO</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t xml:space="preserve"> and it got an overwhelming response from the children.</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00883626" w:rsidP="00A22C3A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>ALARM</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00001D44" w:rsidP="00A22C3A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00A22C3A"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Developed</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00A22C3A"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> an alarm in android with minimalist user interface.</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00883626" w:rsidP="00A22C3A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00883626" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">LIGHT </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>A</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00001D44"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> DREAM</w:t></w:r></w:p><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00A22C3A"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>This</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00A22C3A"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> rails web application helps sports persons to live their dream by registering and getting a trainer, a sponsor and also a view of the events lined up around the world at that very moment. </w:t></w:r><w:r w:rsidR="00582001" w:rsidRPr="00A22C3A"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Presented</w:t></w:r><w:r w:rsidR="0013737E"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> it in </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>CMA</w:t></w:r><w:r w:rsidR="0013737E"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0013737E"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>hackathon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3420" w:type="dxa"/></w:tcPr><w:p w:rsidR="004F17C4" w:rsidRPr="00627399" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>PROFESSIONAL ACHIEVEMENTS</w:t></w:r></w:p><w:p w:rsidR="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="004F17C4"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Placed first rank in Amazon Coding Challenge out of 9K developers in Interview Street</w:t></w:r></w:p><w:p w:rsidR="004F17C4" w:rsidRPr="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:ind w:left="342"/><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00485411" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="004F17C4"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Got ‘A’ band which is highest in rank in TCS for the last year appraisal</w:t></w:r></w:p><w:p w:rsidR="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:ind w:left="342"/><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p><w:p w:rsidR="004F17C4" w:rsidRPr="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="004F17C4"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>Holding ‘Ace of the Quarter – 2013’ award which is given for performing well and good in all the corporate parameters among overall Assurance Services Unit(ASU) level</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00485411" w:rsidRPr="00A22C3A" w:rsidTr="00E779E0"><w:tc><w:tcPr><w:tcW w:w="7380" w:type="dxa"/></w:tcPr><w:p w:rsidR="00A22C3A" w:rsidRPr="00627399" w:rsidRDefault="00001D44" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:bCs/><w:caps/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>SOFT SKILLS</w:t></w:r></w:p><w:p w:rsidR="00485411" w:rsidRPr="00001D44" w:rsidRDefault="00001D44" w:rsidP="00001D44"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>I</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00001D44"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> actively initiate every team outing and celebration, managing the funds of the same thus showing my organizing and management skills.</w:t></w:r></w:p><w:p w:rsidR="00A22C3A" w:rsidRPr="00A22C3A" w:rsidRDefault="00627399" w:rsidP="0013737E"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>I</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00A22C3A"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> carry a lot of humor along in the form of instantaneous one liner</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>s and mimicry of great people. M</w:t></w:r><w:r w:rsidR="00A22C3A" w:rsidRPr="00A22C3A"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">y colleagues appreciate my skills and relish my company. </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3420" w:type="dxa"/></w:tcPr><w:p w:rsidR="004F17C4" w:rsidRPr="004F17C4" w:rsidRDefault="004F17C4" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="004F17C4"><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>CONTACT</w:t></w:r></w:p><w:p w:rsidR="004F17C4" w:rsidRDefault="00627399" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">5/501-2, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Fathima</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> Nagar,</w:t></w:r></w:p><w:p w:rsidR="00627399" w:rsidRPr="00627399" w:rsidRDefault="00627399" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Malaipatti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> Road,</w:t></w:r></w:p><w:p w:rsidR="00627399" w:rsidRPr="00627399" w:rsidRDefault="00627399" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Thottanoothu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> post,</w:t></w:r></w:p><w:p w:rsidR="00627399" w:rsidRPr="00627399" w:rsidRDefault="00627399" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Balakrishna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Puram</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>,</w:t></w:r></w:p><w:p w:rsidR="00485411" w:rsidRPr="000B2CA0" w:rsidRDefault="00627399" w:rsidP="00E779E0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="342"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00627399"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Dindigul – 624005</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00485411" w:rsidRPr="00A22C3A" w:rsidTr="00E779E0"><w:trPr><w:trHeight w:val="810"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="7380" w:type="dxa"/></w:tcPr><w:p w:rsidR="00883626" w:rsidRPr="00A22C3A" w:rsidRDefault="00883626" w:rsidP="00A22C3A"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3420" w:type="dxa"/></w:tcPr><w:p w:rsidR="00485411" w:rsidRPr="00A22C3A" w:rsidRDefault="00485411" w:rsidP="00E779E0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p w:rsidR="005B2A67" w:rsidRPr="00883626" w:rsidRDefault="00AB0D34" w:rsidP="000B2CA0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:sectPr w:rsidR="005B2A67" w:rsidRPr="00883626" w:rsidSect="000B2CA0"><w:footerReference w:type="default" r:id="rId9"/><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="990" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="360"/></w:sectPr>'

$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" ' +
       'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
       'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
       '<pkg:part pkg:name="/word/document.xml" ' +
       'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData><w:document><w:body>' + $newBodyXml + '</w:body></w:document></pkg:xmlData>' +
       '</pkg:part></pkg:package>'

$d.Content.InsertXML($pkg)

# Step 6: top margin 990 -> 630 twips (49.5pt -> 31.5pt), set through the
# regular PageSetup property rather than by hand-editing <w:pgMar> XML.
$d.PageSetup.TopMargin = 31.5

Write-Output "Applied about-me edit; top margin now $($d.PageSetup.TopMargin)pt"
